# 12th commit - Changes in collection agency add new agent and agent list tc
# Adds a new worksheet "CA-AddNewAgent" (mirrors the other "project /
# TestScenario / Run / Name" cover sheets already in the workbook) at the
# end of the workbook, populates it, and leaves it as the active sheet /
# active selection - matching the selection change left on
# CoreHOUserCreation as well.

$wb = $excel.ActiveWorkbook

# --- Update selection on CoreHOUserCreation (loses tabSelected once we
#     activate the new sheet further down, but keeps this selection) -----
$wsCore = $wb.Worksheets.Item("CoreHOUserCreation")
$wsCore.Activate()
$wsCore.Range("F9").Select()

# --- Add the new worksheet at the very end of the workbook --------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "CA-AddNewAgent"

# --- Bring in the header / body formatting used by the sibling sheets ---
# Row 1 (header) formatting matches the CoreUserManagement sheet's header.
$wsUserMgmt = $wb.Worksheets.Item("CoreUserManagement")
$wsUserMgmt.Range("A1:D1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# Row 2 (data) formatting matches the Disposition_master sheet's data row.
$wsDisposition = $wb.Worksheets.Item("Disposition_master")
$wsDisposition.Range("A2:D2").Copy()
$newSheet.Range("A2:D2").PasteSpecial(-4122)

# --- Populate the values --------------------------------------------------
$newSheet.Range("A1").Value = "project"
$newSheet.Range("B1").Value = "TestScenario"
$newSheet.Range("C1").Value = "Run"
$newSheet.Range("D1").Value = "Name"

$newSheet.Range("A2").Value = "Beacon FCM"
$newSheet.Range("B2").Value = "Collection Agency Agent management Add new agent and agent list"
$newSheet.Range("C2").Value = "Yes"
$newSheet.Range("D2").Value = "muthoot4"

# --- Column widths (approximate the sized columns on the sheet) ---------
$newSheet.Columns.Item(1).ColumnWidth = 11.09
$newSheet.Columns.Item(2).ColumnWidth = 60.09
$newSheet.Columns.Item(4).ColumnWidth = 17.36

# --- Make the new sheet the active tab / selection -----------------------
$newSheet.Activate()
$newSheet.Range("E2").Select()
